{"js": "// Update the division-problem table: each data cell's \"a\u00f7b=\" expression is\n// replaced with a new \"a\u00f7b=\" expression, in table (row, column) order.\n// Rows 0, 4, 8, 12, 16 hold the 5 problems per \"page row\"; the other rows\n// are spacer rows and are left untouched.\nconst replacements = [\n  [0, 0, \"55\u00f72=\"],\n  [0, 1, \"12\u00f75=\"],\n  [0, 2, \"85\u00f73=\"],\n  [0, 3, \"31\u00f77=\"],\n  [0, 4, \"66\u00f74=\"],\n  [4, 0, \"92\u00f72=\"],\n  [4, 1, \"35\u00f79=\"],\n  [4, 2, \"36\u00f77=\"],\n  [4, 3, \"18\u00f75=\"],\n  [4, 4, \"69\u00f72=\"],\n  [8, 0, \"61\u00f78=\"],\n  [8, 1, \"90\u00f77=\"],\n  [8, 2, \"50\u00f72=\"],\n  [8, 3, \"25\u00f74=\"],\n  [8, 4, \"51\u00f76=\"],\n  [12, 0, \"73\u00f73=\"],\n  [12, 1, \"72\u00f75=\"],\n  [12, 2, \"71\u00f75=\"],\n  [12, 3, \"59\u00f79=\"],\n  [12, 4, \"19\u00f72=\"],\n  [16, 0, \"46\u00f73=\"],\n  [16, 1, \"25\u00f75=\"],\n  [16, 2, \"91\u00f79=\"],\n  [16, 3, \"52\u00f78=\"],\n  [16, 4, \"19\u00f75=\"],\n];\n\nconst table = context.document.body.tables.items[0];\n\nfor (const [row, col, text] of replacements) {\n  const cell = table.getCell(row, col);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem table: each data cell's \"a\u00f7b=\" expression is\n# replaced with a new \"a\u00f7b=\" expression, addressed by (row, column) using\n# Word's 1-based Table.Cell(row, col) indexing. Rows 1, 5, 9, 13, 17 hold\n# the 5 problems per \"page row\"; the other rows are spacer rows and are\n# left untouched.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @(1, 1, \"55\u00f72=\"),\n    @(1, 2, \"12\u00f75=\"),\n    @(1, 3, \"85\u00f73=\"),\n    @(1, 4, \"31\u00f77=\"),\n    @(1, 5, \"66\u00f74=\"),\n    @(5, 1, \"92\u00f72=\"),\n    @(5, 2, \"35\u00f79=\"),\n    @(5, 3, \"36\u00f77=\"),\n    @(5, 4, \"18\u00f75=\"),\n    @(5, 5, \"69\u00f72=\"),\n    @(9, 1, \"61\u00f78=\"),\n    @(9, 2, \"90\u00f77=\"),\n    @(9, 3, \"50\u00f72=\"),\n    @(9, 4, \"25\u00f74=\"),\n    @(9, 5, \"51\u00f76=\"),\n    @(13, 1, \"73\u00f73=\"),\n    @(13, 2, \"72\u00f75=\"),\n    @(13, 3, \"71\u00f75=\"),\n    @(13, 4, \"59\u00f79=\"),\n    @(13, 5, \"19\u00f72=\"),\n    @(17, 1, \"46\u00f73=\"),\n    @(17, 2, \"25\u00f75=\"),\n    @(17, 3, \"91\u00f79=\"),\n    @(17, 4, \"52\u00f78=\"),\n    @(17, 5, \"19\u00f75=\")\n)\n\nforeach ($item in $replacements) {\n    $row = $item[0]\n    $col = $item[1]\n    $text = $item[2]\n    $t.Cell($row, $col).Range.Text = $text\n}\n"}
